$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 56, shifting existing rows 56:83 down to 57:84
$ws.Rows(56).Insert()

# Populate the newly inserted row 56 with the new record
$ws.Range("A56").Value = 5
$ws.Range("B56").Value = "Macroferia Regional de Talca"
$ws.Range("C56").Value = "Maule"
$ws.Range("D56").Value = 44572
$ws.Range("E56").Value = 7
$ws.Range("F56").Value = 100112001
$ws.Range("G56").Value = "Berenjena"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 200
$ws.Range("K56").Value = 8000
$ws.Range("L56").Value = 8000
$ws.Range("M56").Value = 8000
$ws.Range("N56").Value = "`$/caja 50 unidades"
$ws.Range("O56").Value = "Región del Maule"
$ws.Range("P56").Value = 160
$ws.Range("Q56").Value = 50
$ws.Range("R56").Value = "Hortaliza"
